$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the first-line indent (w:ind w:firstLine="567") from the two
#    "Таблица 1.x" table-caption paragraphs.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $txt = $p.Range.Text
    if ($txt -like "Таблица 1.*" -and [math]::Round($p.Format.FirstLineIndent, 2) -eq 28.35) {
        $p.Format.FirstLineIndent = 0
    }
}

# ---------------------------------------------------------------------------
# 2) Merge the "3.1." run into the preceding run's text (" ...таблице 1." +
#    "3.1." -> " ...таблице 1.3.1.") and drop the now-empty _GoBack bookmark
#    that used to sit between them.
# ---------------------------------------------------------------------------
$insertedText = "3.1."

$gb = $d.Bookmarks("_GoBack")
$mergePoint = $gb.End
$oldRunRange = $d.Range($mergePoint, $mergePoint + $insertedText.Length)
$oldRunText = $oldRunRange.Text

# Insert "3.1." right at the bookmark end (this lands inside/extends the run
# that ends with "...таблице 1."), then drop the bookmark and delete the
# original standalone "3.1." run, whose text now immediately follows.
$ip = $d.Range($mergePoint, $mergePoint)
$ip.InsertAfter($insertedText)

$gb2 = $d.Bookmarks("_GoBack")
$gb2.Delete()

$leftoverStart = $mergePoint + $insertedText.Length
$leftover = $d.Range($leftoverStart, $leftoverStart + $oldRunText.Length)
$leftover.Delete()

# ---------------------------------------------------------------------------
# 3) Re-create the _GoBack bookmark (now collapsed/empty) at the start of the
#    following paragraph (the "Таблица 1.2..." caption paragraph).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Таблица 1.*" -and $p.Range.Text -like "*Описание выходного документа*") {
        $startPoint = $d.Range($p.Range.Start, $p.Range.Start)
        $d.Bookmarks.Add("_GoBack", $startPoint)
    }
}
